$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Sheet2"
$ws2.Move($null, $ws1)

$ws2b = $wb.Worksheets.Item("Sheet2")
$ws2b.Range("A1").Formula = "=FV(-1,-2,1)"
$ws2b.Range("B1").Formula = "=FV(-3, -3.5, 1)"
$ws2b.Range("A1:B1").Style = "Calculation"
